$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Link Budget Analysis: Antenna Gains and Cable losses
# Row 15 corresponds to "exp_pl" (expected path loss, dB)
$ws.Range("C15").Value = 94.08220451406545
$ws.Range("D15").Value = 23.00795847371039
$ws.Range("E15").Value = 45.26
$ws.Range("F15").Value = 77.25999999999999
$ws.Range("G15").Value = 90.25999999999999
$ws.Range("H15").Value = 110.26
$ws.Range("I15").Value = 145.26
